$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "26.307.32"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.588.67"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.72%  "
Set-TextValue $ws.Range("D5") "209.97"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").Value = "  -0.32%  "
Set-TextValue $ws.Range("D10") "19.39"
$ws.Range("E10").Value = "  -0.46%  "
Set-TextValue $ws.Range("D11") "0.0848"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.812.02"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").Value = "1.584.71"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "26.314.01"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  +5.99%  "
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  +1.03%  "
Set-TextValue $ws.Range("D24") "2.14"
$ws.Range("E24").Value = "  -2.56%  "
Set-TextValue $ws.Range("D25") "144.44"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").Value = "1.316.00"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("E35").Value = "  -2.02%  "
Set-TextValue $ws.Range("D36") "0.608"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  -13.00%  "
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("E42").Value = "  +4.49%  "
Set-TextValue $ws.Range("D43") "0.766"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  -0.59%  "
Set-TextValue $ws.Range("D45") "62.24"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "1.724.62"
$ws.Range("E46").Value = "  -0.09%  "
Set-TextValue $ws.Range("D47") "87.35"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("E49").Value = "  -1.00%  "
Set-TextValue $ws.Range("D50") "0.0977"
$ws.Range("E50").Value = "  -4.04%  "
$ws.Range("E51").Value = "  -0.78%  "
